$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Dll1"
$ws.Range("C2").Value = "Notch2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 10.34874566666667
$ws.Range("H2").Value = 31.046237
$ws.Range("I2").Value = 0.614862320492409
$ws.Range("J2").Value = 0.6148623204924089
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 31.96029466666667
$ws.Range("N2").Value = 95.88088399999999
$ws.Range("O2").Value = 0.3907265741426954
$ws.Range("P2").Value = 0.3907265741426953
$ws.Range("Q2").Value = 330.7489609370565
$ws.Range("R2").Value = 2976.740648433508
$ws.Range("S2").Value = 0.240243048055427
$ws.Range("T2").Value = 0.2402430480554269

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Dll1"
$ws.Range("C3").Value = "Notch2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 10.34874566666667
$ws.Range("H3").Value = 31.046237
$ws.Range("I3").Value = 0.614862320492409
$ws.Range("J3").Value = 0.6148623204924089
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 26.34807
$ws.Range("N3").Value = 79.04420999999999
$ws.Range("O3").Value = 0.3221150253382706
$ws.Range("P3").Value = 0.3221150253382706
$ws.Range("Q3").Value = 272.66947523753
$ws.Range("R3").Value = 2454.02527713777
$ws.Range("S3").Value = 0.1980563919449602
$ws.Range("T3").Value = 0.1980563919449601

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Dll1"
$ws.Range("C4").Value = "Notch2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 10.34874566666667
$ws.Range("H4").Value = 31.046237
$ws.Range("I4").Value = 0.614862320492409
$ws.Range("J4").Value = 0.6148623204924089
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 23.48872
$ws.Range("N4").Value = 70.46616
$ws.Range("O4").Value = 0.287158400519034
$ws.Range("P4").Value = 0.287158400519034
$ws.Range("Q4").Value = 243.0787893155467
$ws.Range("R4").Value = 2187.70910383992
$ws.Range("S4").Value = 0.1765628804920218
$ws.Range("T4").Value = 0.1765628804920218

# Row 5
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Dll1"
$ws.Range("C5").Value = "Notch2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 6.482251000000001
$ws.Range("H5").Value = 19.446753
$ws.Range("I5").Value = 0.3851376795075911
$ws.Range("J5").Value = 0.385137679507591
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 31.96029466666667
$ws.Range("N5").Value = 95.88088399999999
$ws.Range("O5").Value = 0.3907265741426954
$ws.Range("P5").Value = 0.3907265741426953
$ws.Range("Q5").Value = 207.1746520632947
$ws.Range("R5").Value = 1864.571868569652
$ws.Range("S5").Value = 0.1504835260872684
$ws.Range("T5").Value = 0.1504835260872684

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Dll1"
$ws.Range("C6").Value = "Notch2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 6.482251000000001
$ws.Range("H6").Value = 19.446753
$ws.Range("I6").Value = 0.3851376795075911
$ws.Range("J6").Value = 0.385137679507591
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 26.34807
$ws.Range("N6").Value = 79.04420999999999
$ws.Range("O6").Value = 0.3221150253382706
$ws.Range("P6").Value = 0.3221150253382706
$ws.Range("Q6").Value = 170.79480310557
$ws.Range("R6").Value = 1537.15322795013
$ws.Range("S6").Value = 0.1240586333933104
$ws.Range("T6").Value = 0.1240586333933104

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Dll1"
$ws.Range("C7").Value = "Notch2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 6.482251000000001
$ws.Range("H7").Value = 19.446753
$ws.Range("I7").Value = 0.3851376795075911
$ws.Range("J7").Value = 0.385137679507591
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 23.48872
$ws.Range("N7").Value = 70.46616
$ws.Range("O7").Value = 0.287158400519034
$ws.Range("P7").Value = 0.287158400519034
$ws.Range("Q7").Value = 152.25977870872
$ws.Range("R7").Value = 1370.33800837848
$ws.Range("S7").Value = 0.1105955200270122
$ws.Range("T7").Value = 0.1105955200270122
